$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for rows 2-5 to reflect repulled data
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = 3
